# Add one more batchsql test case (batch_012), and renumber the
# batch_010/batch_011 "Batch_sql" / "Table_name" / "Query_sql1"
# identifiers to the zero-padded 3-digit scheme (batch_sql_10 ->
# batch_sql_010, batch11 -> batch011, ...) to make room for batch_012.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 13 (batch_012 case), mirroring the batch_011 row -- written
# first (in this column order) so new shared strings land in the same
# slot order as the source workbook.
$ws.Range("A13").Value = "batch_012"
$ws.Range("C13").Value = "批量操作语句12执行"
$ws.Range("K13").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_012.csv"

# Row 12 (batch_011): Table_name / Query_sql1 get the padded "011" id.
$ws.Range("G12").Value = "batch011"

$ws.Range("G13").Value = "batch012"

$ws.Range("J12").Value = "select * from `$batch011"
$ws.Range("J13").Value = "select * from `$batch012"

# Row 11 (batch_010) / 12 / 13: Batch_sql identifiers get padded "0xx".
$ws.Range("I11").Value = "batch_sql_010"
$ws.Range("I12").Value = "batch_sql_011"
$ws.Range("I13").Value = "batch_sql_012"

# Remaining row 13 cells (values reuse existing shared strings).
$ws.Range("B13").Value = "y"
$ws.Range("D13").Value = "batchsql"
$ws.Range("E13").Value = "SingleTable"
$ws.Range("N13").Value = "csv_containsAll"

# Match the saved selection state from the edit.
$ws.Range("F16").Select()
